$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.983.29"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.678.19"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'215.05"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'20.33"
$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "1.913.49"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.676.26"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "'0.529"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "'65.72"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'8.24"
$ws.Range("E17").Value = "  +6.62%  "
$ws.Range("D18").Value = "26.998.82"
$ws.Range("D19").Value = "'235.53"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'4.45"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  -3.30%  "
$ws.Range("D25").Value = "'146.35"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'16.07"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "1.478.00"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +4.87%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").Value = "'0.583"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "'2.31"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "'67.48"
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").Value = "1.818.81"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("E51").Value = "  +0.06%  "
